# The workbook's Sheet2 is a derived view of Sheet1 (each row N's formulas
# point at Sheet1!<col>N). The target edit removes the row whose
# injection_rate was 0.03 (original row 3). Excel keeps every other row's
# formula text untouched when a row is deleted above it, so the rows below
# simply slide up one slot while still referencing their original Sheet1
# row - exactly the "Sheet1!<col>(N+1)" pattern seen after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# Delete the injection_rate = 0.03 row; dimension shrinks from I22 to I21
# and all following rows shift up automatically.
$ws.Rows.Item(3).Delete() | Out-Null

# Match the refreshed view state: zoomed out to 115% with the selection
# resting on J2.
$excel.ActiveWindow.Zoom = 115
$ws.Range("J2").Select() | Out-Null
